# Generate Report for Handoff
# Adds a new file's status row (GUID 0ce84126-...) to the Overview, zh-cn and
# de-de tables/sheets, mirroring the existing e83111b4-... row structure.
# This models "handing off" a newly-ready file for localization.

$wb = $excel.ActiveWorkbook

$mdFileName    = '0ce84126-cc1b-4a4c-a477-5389f3166f1aooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$mdPathAndName = 'e2e\0ce84126-cc1b-4a4c-a477-5389f3166f1aooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$xlfGuidZh     = '0ce84126-cc1b-4a4c-a477-5389f3166f1aoooooooooooooooooooooooooooooooooooooooo.ff4e316acd544ac594455348fa5be8fdc1e9b2bd.zh-cn.xlf'
$xlfGuidDe     = '0ce84126-cc1b-4a4c-a477-5389f3166f1aoooooooooooooooooooooooooooooooooooooooo.ff4e316acd544ac594455348fa5be8fdc1e9b2bd.de-de.xlf'

$hyperlinkTarget = 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b47f4f26d735471242b5bee8117374cd183a5df3/e2e/' + $mdFileName

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)
$tOverview = $wsOverview.ListObjects.Item(1)
$tOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $mdFileName
$wsOverview.Range("B3").Value = $mdPathAndName
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-18 06:27:34"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $hyperlinkTarget, [Type]::Missing, [Type]::Missing, $mdPathAndName) | Out-Null

$wsOverview.Columns.Item(5).AutoFit() | Out-Null
$wsOverview.Columns.Item(6).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item(2)
$tZhCn = $wsZhCn.ListObjects.Item(1)
$tZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = $mdFileName
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = $xlfGuidZh
$wsZhCn.Range("H3").Value = "2016-08-18 06:27:29"
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "False"
$wsZhCn.Range("P3").Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $hyperlinkTarget, [Type]::Missing, [Type]::Missing, $mdFileName) | Out-Null

$wsZhCn.Columns.Item(3).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item(3)
$tDeDe = $wsDeDe.ListObjects.Item(1)
$tDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = $mdFileName
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = $xlfGuidDe
$wsDeDe.Range("H3").Value = "2016-08-18 06:27:34"
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $hyperlinkTarget, [Type]::Missing, [Type]::Missing, $mdFileName) | Out-Null

$wsDeDe.Columns.Item(3).AutoFit() | Out-Null

Write-Output "Report generated for handoff"
